$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The period column (E16:E27) previously listed periods 2302 -> 2401 in
# ascending order. The new database lists the most recent period first,
# i.e. descending order 2401 -> 2302.
$periods = @("2401","2312","2311","2310","2309","2308","2307","2306","2305","2304","2303","2302")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}
